$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.300.17'
$ws.Range("D3").Value = '4.054.53'
$ws.Range("E3").Value = '  +5.57%  '
$ws.Range("E4").Value = '  +0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.92%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.44'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.03%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.710'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +17.84%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").Value = '  +9.24%  '
$ws.Range("E10").Value = '  +7.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000335'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '50.61'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +22.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.06'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +8.16%  '
$ws.Range("D14").Value = '4.700.07'
$ws.Range("E14").Value = '  +5.77%  '
$ws.Range("D15").Value = '4.048.01'
$ws.Range("E15").Value = '  +5.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.42'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '21.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.24'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.49%  '
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("D20").Value = '72.227.22'
$ws.Range("E20").Value = '  +5.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '440.37'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '101.11'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +17.49%  '
$ws.Range("E23").Value = '  +7.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +8.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.25'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +8.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.58'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.82%  '
$ws.Range("E29").Value = '  +3.13%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.37'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +20.46%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.04%  '
$ws.Range("E32").Value = '  +7.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '675.90'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.80'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.80%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '66.78'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.37%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '42.88'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.76%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.442'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.98%  '
$ws.Range("D38").Value = '0.0₃0874'
$ws.Range("E38").Value = '  +6.09%  '
$ws.Range("E39").Value = '  +7.37%  '
$ws.Range("E40").Value = '  +3.95%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0509'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.29%  '
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.16'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.155'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +12.42%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.81'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.90%  '
$ws.Range("E47").Value = '  +2.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.55'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +13.53%  '
$ws.Range("E49").Value = '  +6.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.39'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000275'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.94%  '
